$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Delete row 2 (the "抚州路临次元06路新春国漫嘉年华" entry); all rows below shift up by one
    $ws.Rows.Item(2).Delete()

    # Renumber the index column (A) to stay sequential: A2=1, A3=2, ... A36=35
    for ($r = 2; $r -le 36; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Refresh the "想去人数" (column F) counts that changed between scrapes
    $ws.Cells.Item(2, 6).Value = 55
    $ws.Cells.Item(4, 6).Value = 101
    $ws.Cells.Item(5, 6).Value = 137
    $ws.Cells.Item(6, 6).Value = 1334
    $ws.Cells.Item(7, 6).Value = 1571
    $ws.Cells.Item(9, 6).Value = 440
    $ws.Cells.Item(13, 6).Value = 74
    $ws.Cells.Item(14, 6).Value = 119
    $ws.Cells.Item(15, 6).Value = 284
    $ws.Cells.Item(18, 6).Value = 1769
    $ws.Cells.Item(19, 6).Value = 74
    $ws.Cells.Item(22, 6).Value = 691
    $ws.Cells.Item(25, 6).Value = 4269
    $ws.Cells.Item(27, 6).Value = 289
    $ws.Cells.Item(28, 6).Value = 1122
    $ws.Cells.Item(31, 6).Value = 632
    $ws.Cells.Item(33, 6).Value = 323
    $ws.Cells.Item(35, 6).Value = 160
}
